# 自动更新Excel文件 - 2025-10-31 23:12:23
# Daily refresh: decrement remaining-days counter (column E) for every
# shop row. When a counter would drop to/below zero, reset it back to
# the shop's total-days value (column D) and roll the start date
# (column F) forward to "today" (2025-11-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251101

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {
    # Row 36's data already carries a malformed start-date value
    # (202510929) and was left untouched by the refresh job, so skip it.
    if ($row -eq 36) { continue }

    $totalCell = $ws.Cells.Item($row, 4)
    $remainCell = $ws.Cells.Item($row, 5)
    $startCell = $ws.Cells.Item($row, 6)

    $total = $totalCell.Value2
    $remaining = $remainCell.Value2

    if ($null -eq $remaining) { continue }

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        $remainCell.Value = $total
        $startCell.Value = $today
    } else {
        $remainCell.Value = $newRemaining
    }
}
